{"js": "// Add a new bulleted list item after the last paragraph in the document\n// (\"Object Pool Design Pattern through Queue Data Structure.\"), containing\n// the text \"Separation of data from behaviours using Scriptable Objects\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert a new paragraph after the last paragraph; it inherits the\n// paragraph's list/style formatting (ListParagraph style + numPr bullet).\nconst newParagraph = lastParagraph.insertParagraph(\n  \"Separation of data from behaviours using Scriptable Objects\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Add a new bulleted list item after the last paragraph in the document\n# (\"Object Pool Design Pattern through Queue Data Structure.\"), containing\n# the text \"Separation of data from behaviours using Scriptable Objects\".\n$d = $word.ActiveDocument\n$lastPara = $d.Paragraphs.Last\n\n# Insert a new paragraph mark after the last paragraph; the new paragraph\n# inherits the preceding paragraph's formatting (ListParagraph style +\n# bullet numPr), matching how Word splits a paragraph.\n$lastPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Last\n$newPara.Range.Text = \"Separation of data from behaviours using Scriptable Objects\"\n"}
